$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.456404
$ws.Range("H2").Value = 1.369212
$ws.Range("I2").Value = 0.01914960767004715
$ws.Range("J2").Value = 0.01914960767004715
$ws.Range("M2").Value = 4.959409333333333
$ws.Range("N2").Value = 14.878228
$ws.Range("O2").Value = 0.8271666313262851
$ws.Range("P2").Value = 0.8271666313262852
$ws.Range("Q2").Value = 2.263494257370667
$ws.Range("R2").Value = 20.371448316336
$ws.Range("S2").Value = 0.01583991646765289
$ws.Range("T2").Value = 0.01583991646765289
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.456404
$ws.Range("H3").Value = 1.369212
$ws.Range("I3").Value = 0.01914960767004715
$ws.Range("J3").Value = 0.01914960767004715
$ws.Range("O3").Value = 0.09421438109281059
$ws.Range("P3").Value = 0.09421438109281059
$ws.Range("Q3").Value = 0.2578122744426667
$ws.Range("R3").Value = 2.320310469984
$ws.Range("S3").Value = 0.001804168434803631
$ws.Range("T3").Value = 0.001804168434803631
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.456404
$ws.Range("H4").Value = 1.369212
$ws.Range("I4").Value = 0.01914960767004715
$ws.Range("J4").Value = 0.01914960767004715
$ws.Range("O4").Value = 0.07861898758090437
$ws.Range("P4").Value = 0.07861898758090438
$ws.Range("Q4").Value = 0.2151363705573333
$ws.Range("R4").Value = 1.936227335016
$ws.Range("S4").Value = 0.001505522767590628
$ws.Range("T4").Value = 0.001505522767590628
$ws.Range("I5").Value = 0.8285024587002443
$ws.Range("J5").Value = 0.8285024587002443
$ws.Range("M5").Value = 4.959409333333333
$ws.Range("N5").Value = 14.878228
$ws.Range("O5").Value = 0.8271666313262851
$ws.Range("P5").Value = 0.8271666313262852
$ws.Range("Q5").Value = 97.92945055572847
$ws.Range("R5").Value = 881.3650550015561
$ws.Range("S5").Value = 0.6853095878086257
$ws.Range("T5").Value = 0.6853095878086258
$ws.Range("I6").Value = 0.8285024587002443
$ws.Range("J6").Value = 0.8285024587002443
$ws.Range("O6").Value = 0.09421438109281059
$ws.Range("P6").Value = 0.09421438109281059
$ws.Range("S6").Value = 0.07805684638031539
$ws.Range("T6").Value = 0.07805684638031539
$ws.Range("I7").Value = 0.8285024587002443
$ws.Range("J7").Value = 0.8285024587002443
$ws.Range("O7").Value = 0.07861898758090437
$ws.Range("P7").Value = 0.07861898758090438
$ws.Range("S7").Value = 0.06513602451130324
$ws.Range("T7").Value = 0.06513602451130325
$ws.Range("I8").Value = 0.1523479336297086
$ws.Range("J8").Value = 0.1523479336297086
$ws.Range("M8").Value = 4.959409333333333
$ws.Range("N8").Value = 14.878228
$ws.Range("O8").Value = 0.8271666313262851
$ws.Range("P8").Value = 0.8271666313262852
$ws.Range("Q8").Value = 18.007610329924
$ws.Range("R8").Value = 162.068492969316
$ws.Range("S8").Value = 0.1260171270500066
$ws.Range("T8").Value = 0.1260171270500066
$ws.Range("I9").Value = 0.1523479336297086
$ws.Range("J9").Value = 0.1523479336297086
$ws.Range("O9").Value = 0.09421438109281059
$ws.Range("P9").Value = 0.09421438109281059
$ws.Range("S9").Value = 0.01435336627769158
$ws.Range("T9").Value = 0.01435336627769158
$ws.Range("I10").Value = 0.1523479336297086
$ws.Range("J10").Value = 0.1523479336297086
$ws.Range("O10").Value = 0.07861898758090437
$ws.Range("P10").Value = 0.07861898758090438
$ws.Range("S10").Value = 0.01197744030201051
$ws.Range("T10").Value = 0.01197744030201051
